$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in H6 and I6 (they are removed entirely from the sheet)
$ws.Range("H6:I6").ClearContents()

# Update the selected cell/range shown in the sheet view
$ws.Range("H6").Select()
